# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Column D holds price strings (e.g. "12.30", "0.06120",
    # "27.909.09") that must stay literal text -- left as a plain
    # Value assignment, Excel would auto-convert numeric-looking
    # strings into real numbers and lose trailing zeros / precision.
    # Forcing NumberFormat="@" keeps it text, then restoring the
    # "Normal" style afterwards avoids leaving a stray number format
    # applied to the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.909.09"
$ws.Range("E2").Value = "  +2.05%  "
Set-TextValue $ws.Range("D3") "1.881.96"
$ws.Range("E3").Value = "  +1.69%  "
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue $ws.Range("D5") "332.98"
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("E6").Value = "  +0.06%  "
Set-TextValue $ws.Range("D7") "0.4745"
$ws.Range("E7").Value = "  +6.15%  "
$ws.Range("E8").Value = "  +4.24%  "
Set-TextValue $ws.Range("D9") "48.14"
$ws.Range("E9").Value = "  -0.83%  "
Set-TextValue $ws.Range("D10") "0.08054"
$ws.Range("E10").Value = "  +2.42%  "
$ws.Range("E11").Value = "  +1.81%  "
Set-TextValue $ws.Range("D12") "21.94"
$ws.Range("E12").Value = "  +2.89%  "
Set-TextValue $ws.Range("D13") "1.880.93"
$ws.Range("E13").Value = "  +2.46%  "
Set-TextValue $ws.Range("D14") "5.968"
$ws.Range("E14").Value = "  +2.29%  "
Set-TextValue $ws.Range("D15") "7.199"
$ws.Range("E15").Value = "  +1.48%  "
Set-TextValue $ws.Range("D16") "1.002"
$ws.Range("E16").Value = "  -0.05%  "
Set-TextValue $ws.Range("D17") "0.00001051"
$ws.Range("E17").Value = "  +2.12%  "
Set-TextValue $ws.Range("D18") "87.33"
$ws.Range("E18").Value = "  +1.92%  "
Set-TextValue $ws.Range("D19") "0.06617"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("E21").Value = "  +0.07%  "
Set-TextValue $ws.Range("D22") "28.033.76"
$ws.Range("E22").Value = "  +2.59%  "
Set-TextValue $ws.Range("D23") "5.516"
$ws.Range("E23").Value = "  +1.22%  "
Set-TextValue $ws.Range("D24") "11.07"
$ws.Range("E24").Value = "  +3.00%  "
Set-TextValue $ws.Range("D25") "2.316"
$ws.Range("E25").Value = "  +2.85%  "
Set-TextValue $ws.Range("D26") "2.130.23"
$ws.Range("E26").Value = "  +3.58%  "
Set-TextValue $ws.Range("D27") "157.56"
$ws.Range("E27").Value = "  +3.93%  "
Set-TextValue $ws.Range("D28") "20.22"
$ws.Range("E28").Value = "  +4.54%  "
Set-TextValue $ws.Range("D29") "2.111"
$ws.Range("E29").Value = "  +2.87%  "
Set-TextValue $ws.Range("D30") "5.626"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("E31").Value = "  +2.95%  "
Set-TextValue $ws.Range("D32") "0.9834"
$ws.Range("E32").Value = "  +5.93%  "
Set-TextValue $ws.Range("D33") "0.09572"
$ws.Range("E33").Value = "  +2.78%  "
Set-TextValue $ws.Range("D34") "1.465"
$ws.Range("E34").Value = "  +0.51%  "
Set-TextValue $ws.Range("D35") "3.631"
$ws.Range("E35").Value = "  +1.04%  "
Set-TextValue $ws.Range("D36") "5.335"
$ws.Range("E36").Value = "  +1.98%  "
Set-TextValue $ws.Range("D37") "0.06120"
$ws.Range("E37").Value = "  +3.00%  "
Set-TextValue $ws.Range("D38") "0.02266"
$ws.Range("E38").Value = "  +2.43%  "
Set-TextValue $ws.Range("D39") "1.234"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("E40").Value = "  -0.43%  "
Set-TextValue $ws.Range("D41") "0.6040"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  +0.05%  "
Set-TextValue $ws.Range("D43") "0.1905"
$ws.Range("E43").Value = "  +3.00%  "
Set-TextValue $ws.Range("D44") "10.33"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D45") "1.260"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D46") "0.5723"
$ws.Range("E46").Value = "  +1.55%  "
Set-TextValue $ws.Range("D47") "12.30"
$ws.Range("E47").Value = "  +1.53%  "
Set-TextValue $ws.Range("D48") "3.418"
$ws.Range("E48").Value = "  +1.73%  "
Set-TextValue $ws.Range("D49") "1.953"
$ws.Range("E49").Value = "  +1.83%  "
Set-TextValue $ws.Range("D50") "0.06835"
$ws.Range("E50").Value = "  -0.21%  "
Set-TextValue $ws.Range("D51") "114.05"
$ws.Range("E51").Value = "  +5.40%  "
